# Generate Report for Handoff
#
# Re-running the localization-status report generation refreshed the rows
# that were still "Ready for handoff" (source rows 4-7) on both locale
# sheets: their Priority is recomputed as "ht", and their Latest Handoff
# Datetime is bumped to the new generation run's timestamp. The Overview
# sheet's "Latest HO Xliff Generate Date" column mirrors the de-de Latest
# Handoff Datetime for the same files, so it moves in lockstep.

$wb = $excel.ActiveWorkbook

$ws_zh = $wb.Worksheets.Item("zh-cn")
$ws_zh.Range("E4:E7").Value = "ht"
$ws_zh.Range("H4:H7").Value = "2016-08-19 06:31:26"

$ws_de = $wb.Worksheets.Item("de-de")
$ws_de.Range("E4:E7").Value = "ht"
$ws_de.Range("H4:H7").Value = "2016-08-19 06:31:31"

$ws_ov = $wb.Worksheets.Item("Overview")
$ws_ov.Range("G4:G7").Value = "2016-08-19 06:31:31"
